$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for changed rows to match the refreshed
# cryptos snapshot. Price-column updates that are plain decimals need the cell coerced
# to Text first so Excel does not reinterpret them (and drop trailing zeros) as numbers -
# the source data keeps these as text (note the thousands-dotted prices like "36.284.27").
$ws.Range("D2").Value = "36.284.27"
$ws.Range("E2").Value = "  -3.42%  "
$ws.Range("D3").Value = "1.967.59"
$ws.Range("E3").Value = "  -3.91%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.78"
$ws.Range("E5").Value = "  -3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  -5.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.30"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.374"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.79"
$ws.Range("E10").Value = "  -4.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0805"
$ws.Range("E11").Value = "  +6.14%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.858"
$ws.Range("E13").Value = "  -6.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.49"
$ws.Range("E14").Value = "  +10.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.09"
$ws.Range("E15").Value = "  -7.00%  "
$ws.Range("D16").Value = "2.255.10"
$ws.Range("E16").Value = "  -3.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.43"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").Value = "1.964.27"
$ws.Range("E18").Value = "  -4.16%  "
$ws.Range("D19").Value = "36.114.76"
$ws.Range("E19").Value = "  -3.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.02"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("D21").Value = "0.0₃0859"
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.38"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.22"
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("E25").Value = "  -4.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  -4.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.84"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.31"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.85"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  +14.35%  "
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.87"
$ws.Range("E32").Value = "  -6.94%  "
$ws.Range("E33").Value = "  -5.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0621"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.40"
$ws.Range("E35").Value = "  -7.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.30"
$ws.Range("E36").Value = "  +6.00%  "
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.10"
$ws.Range("E40").Value = "  +10.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0988"
$ws.Range("E41").Value = "  -3.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.23"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.87"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0213"
$ws.Range("E44").Value = "  -2.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.09"
$ws.Range("E45").Value = "  -4.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.52"
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.03"
$ws.Range("E47").Value = "  -5.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.49"
$ws.Range("E48").Value = "  -7.01%  "
$ws.Range("D49").Value = "1.338.47"
$ws.Range("E49").Value = "  -5.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.84"
$ws.Range("E50").Value = "  -3.79%  "
$ws.Range("D51").Value = "2.149.03"
$ws.Range("E51").Value = "  -3.77%  "

# Rows 37 and 38 swapped: BinanceUSD <-> LidoDAOToken (with refreshed price/volume)
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  -5.95%  "

$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.01%  "
